# TestData.xlsx - add new "Schedule Dashboard via email" test-case rows
# (new EmailFormat options html/xlsx/csv rows, and a new "duplicate job
# name" scheduling sub-section) below the existing TC018 block, mirroring
# the existing layout/styling conventions used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# --- Clone existing row formatting onto the new rows -----------------
# Rows 49:50 are plain "Key / Value" data rows (style used throughout the
# sheet for simple two-column entries) - reuse that formatting for the
# three new EmailFormat rows (51:53).
$ws.Range("A49:B50").Copy() | Out-Null
$ws.Range("A51:B53").PasteSpecial(-4122) | Out-Null

# Row 48 is a merged section-header row ("Send Dashboard via email") -
# reuse that formatting for the new "Schedule Dashboard via email"
# section header at row 54.
$ws.Range("A48:O48").Copy() | Out-Null
$ws.Range("A54:O54").PasteSpecial(-4122) | Out-Null

# More plain "Key / Value" rows for JobName / Email / error-message.
$ws.Range("A49:B50").Copy() | Out-Null
$ws.Range("A55:B57").PasteSpecial(-4122) | Out-Null

# Trailing blank styled cell.
$ws.Range("B49").Copy() | Out-Null
$ws.Range("B58").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Section-header row (54) uses the compact 13.5pt height like the other
# section headers (e.g. row 48); the data rows keep the sheet's normal
# 15.75pt height.
$ws.Rows.Item(54).RowHeight = 13.5

# --- Values ------------------------------------------------------------
$ws.Range("A51").Value = "html"
$ws.Range("B51").Value = "The dashboard will be sent in the same layout it appears now. Insights can be downloaded as images."

$ws.Range("A52").Value = "xlsx"
$ws.Range("B52").Value = "An XLSX file will be sent only for the tables and pivot tables in the dashboard."

$ws.Range("A53").Value = "csv"
$ws.Range("B53").Value = "A CSV file will be sent only for the first table in this dashboard."

$ws.Range("A54").Value = "Schedule Dashboard via email"

$ws.Range("A55").Value = "JobName"
$ws.Range("B55").Value = "Duplicate Job Name"

$ws.Range("A56").Value = "Email"
$ws.Range("B56").Value = "Ahmed.Abdelsalam@incorta.com"

$ws.Range("A57").Value = "DuplicateJobNameErrorMessage"
$ws.Range("B57").Value = "INC_004010050:Another SCHEDULER with the same name [Duplicate Job Name] already exists."

# --- Merge the new section-header row like the other section headers --
$ws.Range("A54:C54").Merge() | Out-Null
